$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update the SWOT "Forças" entry: the tech stack mention moved from
# Python/MySQL to Javascript/Next.js/Puppeteer.
$ws.Range("A4").Value = "Uso de tecnologias consolidadas (Javascript, Next.js, Puppeteer)"

# Move the active selection, matching the saved session state.
$ws.Range("B12").Select()
